$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EVALUATEC")

# Unhide the previously hidden data rows (2-45) that were filtered out by
# the "7 de agosto" date filter -- the filter is being cleared so every
# applicant row becomes visible again.
$ws.Range("A2:A45").EntireRow.Hidden = $false

# Drop the date-group filter on the "FECHA APLICACIÓN OFICIAL" column so the
# table shows every row (no rows excluded by AutoFilter criteria anymore).
$lo = $ws.ListObjects.Item("Tabla2")
$lo.AutoFilter.ShowAllData() | Out-Null

# Move the active selection to match where the author left the cursor.
$ws.Range("G11").Select() | Out-Null
